$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename sheet: drop the "2023-2013" suffix to avoid future errors
$ws.Name = "g3.11c variação"

# Add new column D "Ano" with the year-range label for each data row.
# Give D1 the same header formatting as the existing headers (A1:C1) by
# copying C1's format onto it, rather than creating a brand-new style.
$ws.Range("D1").Value = "Ano"
$ws.Range("C1").Copy()
$ws.Range("D1").PasteSpecial(-4122)

$ws.Range("D2:D10").Value = "2023-2013"
